$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.00258576107567661, 0.997414238924323, 0.991897948629547, 0.990863644199276, 0.0691260127564213, 0.951732459920703, 0.030511980692984, 0.993449405274953, 0.16031718669195, 0.000172384071711774, 0.000517152215135321, 0.000861920358558869, 0, 0.998965695569729, 0.0210308567488364, 0.00379244957765902, 0.00120668850198242, 0.0375797276331667, 0.000172384071711774, 0, 0.023616617824513, 0.00896397172901224, 0.994311325633511)
    3 = @(0, 0, 0.000172384071711774, 0.00655059472504741, 0, 0.00120668850198242, 0, 0.0032752973625237, 0.760213756248923, 0.000861920358558869, 0.00120668850198242, 0.00120668850198242, 0.999310463713153, 0, 0.0543009825892088, 0.0349939665574901, 0.000344768143423548, 0.000172384071711774, 0, 0.000172384071711774, 0.116876400620583, 0.990174107912429, 0)
    4 = @(0.996724702637476, 0.00241337700396483, 0.00758489915531805, 0.00103430443027064, 0.929150146526461, 0.0456817790036201, 0.968970867091881, 0.00275814514738838, 0.00534390622306499, 0, 0.997759007067747, 0.997586622996035, 0, 0.000861920358558869, 0.0239613859679366, 0, 0.000172384071711774, 0.959834511291157, 0.999827615928288, 0.999482847784865, 0.0172384071711774, 0.000172384071711774, 0.00465436993621789)
    5 = @(0, 0.000172384071711774, 0, 0.00120668850198242, 0.000172384071711774, 0.000344768143423548, 0, 0.000172384071711774, 0.0668850198241683, 0.998793311498018, 0.000517152215135321, 0, 0.000517152215135321, 0, 0.896914325116359, 0.959662127219445, 0.99810377521117, 0.000517152215135321, 0, 0, 0.83968281330805, 0.000344768143423548, 0.000689536286847095)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$wb.Save()
